# Refresh FFXIV Leve market-profit figures (columns H-N) across the
# ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets with the latest scraped market-board
# averages. Mirrors the scheduled-runner data refresh described in the
# commit; GSM has no changes this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising
$ws.Range("H2").Value = 3839.8
$ws.Range("J2").Value = 3839.8
$ws.Range("L2").Value = 3839.8
$ws.Range("N2").Value = -4065.8
# Row 7: The Bleat Is On
$ws.Range("H7").Value = 30001
$ws.Range("J7").Value = 30001
$ws.Range("L7").Value = 30001
$ws.Range("N7").Value = -30225
# Row 8: On the Drip
$ws.Range("H8").Value = 16.142857
$ws.Range("I8").Value = 17.666666
$ws.Range("K8").Value = 52.999998
$ws.Range("M8").Value = 86.00000199999999
# Row 12: Don't Be So Tallow
$ws.Range("H12").Value = 316.2
$ws.Range("I12").Value = 292
$ws.Range("J12").Value = 352.5
$ws.Range("K12").Value = 292
$ws.Range("L12").Value = 352.5
$ws.Range("M12").Value = -122
$ws.Range("N12").Value = -692.5
# Row 14: Wand-full Tonight
$ws.Range("H14").Value = 30001
$ws.Range("J14").Value = 30001
$ws.Range("L14").Value = 30001
$ws.Range("N14").Value = -30383
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 596.73334
$ws.Range("I15").Value = 596.73334
$ws.Range("K15").Value = 1790.20002
$ws.Range("M15").Value = -1621.20002
# Row 16: Using Your Arcane Powers for Fun and Profit
$ws.Range("H16").Value = 6333.3335
$ws.Range("J16").Value = 13000
$ws.Range("L16").Value = 13000
$ws.Range("N16").Value = -13460
# Row 20: Shut Up and Take My Gil
$ws.Range("H20").Value = 227.5
$ws.Range("I20").Value = 227.5
$ws.Range("K20").Value = 227.5
$ws.Range("M20").Value = 2.5
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 138
$ws.Range("I33").Value = 119.833336
$ws.Range("K33").Value = 119.833336
$ws.Range("M33").Value = 109.166664
# Row 34: Sophomore Slump
$ws.Range("H34").Value = 5044
$ws.Range("I34").Value = 5044
$ws.Range("K34").Value = 5044
$ws.Range("M34").Value = -4841
# Row 35: Conspicuous Conjuration
$ws.Range("H35").Value = 227.5
$ws.Range("I35").Value = 227.5
$ws.Range("K35").Value = 227.5
$ws.Range("M35").Value = 151.5
# Row 36: You Put Your Left Hand In
$ws.Range("H36").Value = 5044
$ws.Range("I36").Value = 5044
$ws.Range("K36").Value = 5044
$ws.Range("M36").Value = -4329
# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 465.75
$ws.Range("J38").Value = 1648
$ws.Range("L38").Value = 4944
$ws.Range("N38").Value = -5688
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 2742
$ws.Range("I43").Value = 1487
$ws.Range("J43").Value = 3997
$ws.Range("K43").Value = 1487
$ws.Range("L43").Value = 3997
$ws.Range("M43").Value = -1418
$ws.Range("N43").Value = -4135
# Row 54: Arcane Arts for Dummies
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 3957.1875
$ws.Range("I62").Value = 3966.818
$ws.Range("K62").Value = 3966.818
$ws.Range("M62").Value = -3342.818
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 3957.1875
$ws.Range("I65").Value = 3966.818
$ws.Range("K65").Value = 19834.09
$ws.Range("M65").Value = -16714.09
# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 486.16666
$ws.Range("I80").Value = 675.5714
$ws.Range("J80").Value = 221
$ws.Range("K80").Value = 2026.7142
$ws.Range("L80").Value = 663
$ws.Range("M80").Value = -1028.7142
$ws.Range("N80").Value = -2659
# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 486.16666
$ws.Range("I83").Value = 675.5714
$ws.Range("J83").Value = 221
$ws.Range("K83").Value = 6080.1426
$ws.Range("L83").Value = 1989
$ws.Range("M83").Value = -1088.1426
$ws.Range("N83").Value = -11973
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2308.625
$ws.Range("I132").Value = 2191.6086
$ws.Range("K132").Value = 6574.825800000001
$ws.Range("M132").Value = -4044.825800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2544.7144
$ws.Range("I32").Value = 2544.7144
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2544.7144
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2257.7144
$ws.Range("N32").Value = $null
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2414.3333
$ws.Range("I61").Value = 2414.3333
$ws.Range("K61").Value = 2414.3333
$ws.Range("M61").Value = -2202.3333
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2782.6667
$ws.Range("I74").Value = 2782.6667
$ws.Range("K74").Value = 2782.6667
$ws.Range("M74").Value = -1908.6667
# Row 76: Sometimes the South Wins
$ws.Range("H76").Value = 43225.668
$ws.Range("J76").Value = 43225.668
$ws.Range("L76").Value = 43225.668
$ws.Range("N76").Value = -43901.668
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2782.6667
$ws.Range("I77").Value = 2782.6667
$ws.Range("K77").Value = 13913.3335
$ws.Range("M77").Value = -9545.333500000001
# Row 79: The Thriller of Autumn (L)
$ws.Range("H79").Value = 43225.668
$ws.Range("J79").Value = 43225.668
$ws.Range("L79").Value = 43225.668
$ws.Range("N79").Value = -45565.668
# Row 80: A Squire to Inspire
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
# Row 82: Belle of the Brawl
$ws.Range("H82").Value = 34444
$ws.Range("J82").Value = 34444
$ws.Range("L82").Value = 34444
$ws.Range("N82").Value = -35166
# Row 83: All's Fair in Highborn Assassination (L)
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
# Row 85: Shouldering the Shut-ins (L)
$ws.Range("H85").Value = 34444
$ws.Range("J85").Value = 34444
$ws.Range("L85").Value = 34444
$ws.Range("N85").Value = -36940
# Row 97: Ore for Me
$ws.Range("H97").Value = 1177.5714
$ws.Range("I97").Value = 1009.6
$ws.Range("J97").Value = 1597.5
$ws.Range("K97").Value = 1009.6
$ws.Range("L97").Value = 1597.5
$ws.Range("M97").Value = -513.6
$ws.Range("N97").Value = -2589.5
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1293.9131
$ws.Range("I132").Value = 1298.1428
$ws.Range("J132").Value = 1249.5
$ws.Range("K132").Value = 3894.4284
$ws.Range("L132").Value = 3748.5
$ws.Range("M132").Value = -1364.4284
$ws.Range("N132").Value = -8808.5
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2414.3333
$ws.Range("I136").Value = 2414.3333
$ws.Range("K136").Value = 7242.999899999999
$ws.Range("M136").Value = -4692.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 135: Axes to the Maxes
$ws.Range("H135").Value = 47926
$ws.Range("J135").Value = 47926
$ws.Range("L135").Value = 47926
$ws.Range("N135").Value = -58066

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2215.0286
$ws.Range("I31").Value = 1614
$ws.Range("J31").Value = 3526.3635
$ws.Range("K31").Value = 1614
$ws.Range("L31").Value = 3526.3635
$ws.Range("M31").Value = -1319
$ws.Range("N31").Value = -4116.363499999999
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2215.0286
$ws.Range("I34").Value = 1614
$ws.Range("J34").Value = 3526.3635
$ws.Range("K34").Value = 1614
$ws.Range("L34").Value = 3526.3635
$ws.Range("M34").Value = -1412
$ws.Range("N34").Value = -3930.3635
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 4408
$ws.Range("I58").Value = 3885.9048
$ws.Range("J58").Value = 9890
$ws.Range("K58").Value = 3885.9048
$ws.Range("L58").Value = 9890
$ws.Range("M58").Value = -3682.9048
$ws.Range("N58").Value = -10296
# Row 86: Birch, Please
$ws.Range("H86").Value = 9999.700000000001
$ws.Range("I86").Value = 9263.625
$ws.Range("K86").Value = 9263.625
$ws.Range("M86").Value = -8140.625
# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 9999.700000000001
$ws.Range("I89").Value = 9263.625
$ws.Range("K89").Value = 46318.125
$ws.Range("M89").Value = -40702.125
# Row 107: Built to Last
$ws.Range("H107").Value = 1222
$ws.Range("I107").Value = 874.5
$ws.Range("K107").Value = 874.5
$ws.Range("M107").Value = 1045.5
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3727.111
$ws.Range("I134").Value = 3930.5
$ws.Range("J134").Value = 3564.4
$ws.Range("K134").Value = 11791.5
$ws.Range("L134").Value = 10693.2
$ws.Range("M134").Value = -9256.5
$ws.Range("N134").Value = -15763.2
# Row 136: Turali Quality
$ws.Range("H136").Value = 4408
$ws.Range("I136").Value = 3885.9048
$ws.Range("J136").Value = 9890
$ws.Range("K136").Value = 11657.7144
$ws.Range("L136").Value = 29670
$ws.Range("M136").Value = -9107.714399999999
$ws.Range("N136").Value = -34770

$ws = $wb.Worksheets.Item("CUL")
# Row 41: Gegeruju Gets Down
$ws.Range("H41").Value = 2
$ws.Range("I41").Value = 2
$ws.Range("K41").Value = 6
$ws.Range("M41").Value = 332
# Row 44: No More Dumpster Diving
$ws.Range("H44").Value = 457.5
$ws.Range("I44").Value = 360
$ws.Range("J44").Value = 750
$ws.Range("K44").Value = 1080
$ws.Range("L44").Value = 2250
$ws.Range("M44").Value = -682
$ws.Range("N44").Value = -3046
# Row 109: Cure for What Ails
$ws.Range("H109").Value = 1835.5714
$ws.Range("I109").Value = 557.8
$ws.Range("J109").Value = 5030
$ws.Range("K109").Value = 1673.4
$ws.Range("L109").Value = 15090
$ws.Range("M109").Value = -633.3999999999999
$ws.Range("N109").Value = -17170

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1420.2
$ws.Range("J22").Value = 1551
$ws.Range("L22").Value = 1551
$ws.Range("N22").Value = -2141
# Row 27: Fire and Hide
$ws.Range("H27").Value = 1420.2
$ws.Range("J27").Value = 1551
$ws.Range("L27").Value = 1551
$ws.Range("N27").Value = -1765
# Row 44: The Righteous Tools for the Job
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null
# Row 119: Fit for a Friend
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1340.5555
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 3000
$ws.Range("N132").Value = -8060

